$d = $word.ActiveDocument

# 1. Remove the original "_GoBack" bookmark (it currently starts at the very
#    beginning of the document title and ends at the end of the final
#    paragraph). Deleting it removes both the bookmarkStart and bookmarkEnd.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Locate the target paragraph (last paragraph in the document) that
#    currently holds the three runs describing ancient/modern philosophers.
$count = $d.Paragraphs.Count
$target = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*implied differences*") {
        $target = $p
        break
    }
}

if ($null -eq $target) {
    throw "Could not locate the target paragraph (containing 'implied differences')."
}

$r = $target.Range
$paraStart = $r.Start
$paraEnd = $r.End

# 3. Clear the existing text in the paragraph (but keep the paragraph mark).
$clearRange = $d.Range($paraStart, $paraEnd - 1)
$clearRange.Text = ""

# 4. Insert the first new run's text.
$text1 = 'I’m not sure the famous ship example is a good one for Socrates. When he is saying that a good captain uses a high level of knowledge to navigate he is arguing for a scientist ruler, '
$ins1 = $d.Range($paraStart, $paraStart)
$ins1.InsertAfter($text1)

# 5. Insert the second new run's text right after the first.
$junction = $paraStart + $text1.Length
$text2 = 'not a philosopher. Philosophers tend to use thought experiments (like the article) instead of empirical evidence (winds, starts, etc). A philosopher would prefer to write a long winded logical argument about why the ship moves the way it does while a scientist would prefer to understand the mechanics behind the boats movement. '
$ins2 = $d.Range($junction, $junction)
$ins2.InsertAfter($text2)

# 6. Re-create the "_GoBack" bookmark as a zero-length bookmark sitting at the
#    junction between the two runs (this is where Word leaves it after the
#    last edit).
$bmRange = $d.Range($junction, $junction)
$d.Bookmarks.Add("_GoBack", $bmRange)
